$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Update the "Ready for handoff" -> "Handed back: in sync with en-US"
#    status text everywhere it appears (Overview + both language sheets).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (F) / "Latest Handback File" (G) columns
#    on the zh-cn sheet, and link them like the existing handoff columns.
# ---------------------------------------------------------------------------
$mdName      = "be29b3c6-59aa-46b0-b96a-66e4b8108411.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/d14efd809cdcce8d04fea9d1bc91a1a2b6c7d893/e2e/be29b3c6-59aa-46b0-b96a-66e4b8108411.md"
$md2Name     = "ffff444ff36f-7535-4522-bf42-09a0f3962046.md"
$md2Url      = "https://github.com/OpenLocalizationTest/oltest/blob/d14efd809cdcce8d04fea9d1bc91a1a2b6c7d893/e2e/ffff444ff36f-7535-4522-bf42-09a0f3962046.md"

$zhCnXlfName = "be29b3c6-59aa-46b0-b96a-66e4b8108411.a9a21e6b943661e0065b3c42c2b032963dc528ce.zh-cn.xlf"
$zhCnXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/84f16ac585706cf6eaf4df96da10cebd025af5f7/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/be29b3c6-59aa-46b0-b96a-66e4b8108411.a9a21e6b943661e0065b3c42c2b032963dc528ce.zh-cn.xlf"

$deDeXlfName = "be29b3c6-59aa-46b0-b96a-66e4b8108411.a9a21e6b943661e0065b3c42c2b032963dc528ce.de-de.xlf"
$deDeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbd6c6ae57cd70a18363cec5885642079f60c5ec/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/be29b3c6-59aa-46b0-b96a-66e4b8108411.a9a21e6b943661e0065b3c42c2b032963dc528ce.de-de.xlf"

# Latest Handback DateTime (H) on zh-cn keeps its cell reference but the
# underlying text changes; on de-de it becomes a brand-new timestamp.
$wsZhCn.Range("H2").Value = "2016-03-21 05:02:39"
$wsZhCn.Range("H3").Value = "2016-03-21 05:02:39"

# Drop every existing hyperlink on the sheet so the whole collection can be
# rebuilt in the desired final order (this keeps relationship ids stable
# and sequential, matching how Excel numbers them on a fresh write).
$wsZhCn.Range("A1").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhCnXlfUrl, "", "", $zhCnXlfName)

$wsZhCn.Range("F2").Value = $mdName
$wsZhCn.Range("G2").Value = $zhCnXlfName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $mdUrl, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhCnXlfUrl, "", "", $zhCnXlfName)

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $md2Url, "", "", $md2Name)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhCnXlfUrl, "", "", $zhCnXlfName)

$wsZhCn.Range("F3").Value = $mdName
$wsZhCn.Range("G3").Value = $zhCnXlfName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $mdUrl, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhCnXlfUrl, "", "", $zhCnXlfName)

# ---------------------------------------------------------------------------
# 3. Same treatment for the de-de sheet, with its own handback timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Range("H2").Value = "2016-03-21 05:02:53"
$wsDeDe.Range("H3").Value = "2016-03-21 05:02:53"

$wsDeDe.Range("A1").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deDeXlfUrl, "", "", $deDeXlfName)

$wsDeDe.Range("F2").Value = $mdName
$wsDeDe.Range("G2").Value = $deDeXlfName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deDeXlfUrl, "", "", $deDeXlfName)

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $md2Url, "", "", $md2Name)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deDeXlfUrl, "", "", $deDeXlfName)

$wsDeDe.Range("F3").Value = $mdName
$wsDeDe.Range("G3").Value = $deDeXlfName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deDeXlfUrl, "", "", $deDeXlfName)
